$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for price cells whose new value would otherwise be
# auto-parsed as a number by Excel (these are text cells in the source file).
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume change, and the three
# row swaps where the ranking order changed: Dai/RenderToken,
# Filecoin/dogwifhat, Cosmos/OKB/Bittensor).
$ws.Range("D2").Value = "64.661.65"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "3.125.70"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "590.77"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "152.82"
$ws.Range("E6").Value = "  +4.55%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.120.98"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "5.95"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "38.07"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "3.640.89"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "7.26"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "64.210.19"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "3.127.40"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "473.32"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").Value = "14.99"
$ws.Range("E21").Value = "  +4.94%  "
$ws.Range("D22").Value = "0.742"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  +3.06%  "
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  +8.66%  "
$ws.Range("D25").Value = "13.35"
$ws.Range("E25").Value = "  +3.41%  "
$ws.Range("D26").Value = "82.02"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "9.95"
$ws.Range("E28").Value = "  +7.17%  "
$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D30").Value = "2.71"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "0.118"
$ws.Range("E33").Value = "  +8.68%  "
$ws.Range("D34").Value = "27.58"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("D35").Value = "0.0₃0864"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "6.19"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.37"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("D39").Value = "2.30"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "9.37"
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "51.06"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "456.42"
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("D43").Value = "0.297"
$ws.Range("E43").Value = "  +6.18%  "
$ws.Range("D44").Value = "0.0374"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").Value = "2.875.84"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "0.110"
$ws.Range("E46").Value = "  +2.40%  "
$ws.Range("D47").Value = "39.51"
$ws.Range("E47").Value = "  +5.27%  "
$ws.Range("D48").Value = "130.74"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").Value = "25.52"
$ws.Range("E49").Value = "  +5.92%  "
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  +6.06%  "
$ws.Range("E51").Value = "  +0.03%  "
